# Actualizacion automatica del tracker
# - Resolves two previously pending picks (rows 76 and 80) as "Fallo" (-1 profit)
# - Appends the newest tracked pick as a new row (row 81), result still pending

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 76: mark result as Fallo with profit -1
$ws.Cells.Item(76, 7).Value = "Fallo"
$ws.Cells.Item(76, 8).Value = -1

# Row 80: mark result as Fallo with profit -1
$ws.Cells.Item(80, 7).Value = "Fallo"
$ws.Cells.Item(80, 8).Value = -1

# New row 81 with the latest pick (resultado/profit left blank, pending)
$newRow = 81
$ws.Cells.Item($newRow, 1).Value = 14601568
$ws.Cells.Item($newRow, 3).Value = "Rebecca Marino"
$ws.Cells.Item($newRow, 4).Value = "Tatjana Maria"
$ws.Cells.Item($newRow, 5).Value = "Gana Tatjana Maria"
$ws.Cells.Item($newRow, 6).Value = 1.8

# Write the "fecha" column as literal text "2025-09-10" (not an Excel date
# serial). Assigning the string directly makes Excel auto-detect it as a
# date, so instead build it via a formula and paste back the computed
# value, which stores it as plain text - matching every other cell in
# this column.
$ws.Cells.Item($newRow, 2).Formula = '="2025-09-10"'
$ws.Cells.Item($newRow, 2).Copy()
$ws.Cells.Item($newRow, 2).PasteSpecial(-4163)
$excel.CutCopyMode = 0

$wb.Save()
